$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the used cell range, keep column formatting (cols) intact.
$ws.UsedRange.Clear() | Out-Null

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = "#ITEM"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Designator"
$ws.Range("D1").Value = "Manufacturer"
$ws.Range("E1").Value = "Manufacturer Part Number#"
$ws.Range("F1").Value = "Supplier"
$ws.Range("G1").Value = "Supplier Part #"
$ws.Range("H1").Value = "Link"
$ws.Range("I1").Value = "QTY/BOARD"
$ws.Range("J1").Value = "Order QTY"
$ws.Range("K1").Value = "Price each @ 50"
$ws.Range("L1").Value = "price tota"
$ws.Range("M1").Value = "number of boards:"
$ws.Range("N1").Value = 50

# ---- Row 2 ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CAP CER 0.1UF 6.3V 10% X7R 0402"
$ws.Range("C2").Value = "C1, C2, C3"
$ws.Range("D2").Value = "Samsung Electro-Mechanics"
$ws.Range("E2").Value = "CL05B104KQ5NNNC"
$ws.Range("F2").Value = "Digi-Key"
$ws.Range("G2").Value = "1276-1511-1-ND"
$ws.Range("I2").Value = 3
$ws.Range("J2").Formula = "=I2*`$N`$1"
$ws.Range("P2").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# ---- Row 3 ----
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "IMU ACCEL/GYRO/COMPI2C/SPI 24QFN"
$ws.Range("C3").Value = "U9AB1"
$ws.Range("D3").Value = "TDK InvenSense"
$ws.Range("E3").Value = "ICM-20948"
$ws.Range("F3").Value = "Digi-Key"
$ws.Range("G3").Value = "1428-1123-1-ND"
$ws.Range("H3").Value = "https://www.digikey.com/en/products/detail/tdk-invensense/ICM-20948/6623535"
$ws.Range("I3").Value = 1
$ws.Range("J3").Formula = "=I3*`$N`$1"
$ws.Range("P3").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# ---- Row 4 (new row) ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "CAP CER 100UF 6.3V X5R 0805"
$ws.Range("C4").Value = "C_1.1_OUT1, C_1.8_OUT2"
$ws.Range("D4").Value = "Murata Electronics"
$ws.Range("E4").Value = "GRM21BR60J107ME15L"
$ws.Range("F4").Value = "Digi-Key"
$ws.Range("H4").Value = "https://www.digikey.com/en/products/detail/murata-electronics/GRM21BR60J107ME15L/6155751"
$ws.Range("I4").Value = 2
$ws.Range("J4").Formula = "=I4*`$N`$1"
$ws.Range("K4").Value = 0.8974
$ws.Range("L4").Formula = "=K4*`$N`$1"

# ---- Column widths ----
$ws.Columns.Item(3).ColumnWidth = 22.43
$ws.Columns.Item(8).ColumnWidth = 19.6

# ---- View settings ----
$ws.Range("B4").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
